$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="35.041.36"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '="  +0.46%  "'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("D3").Formula = '="1.848.68"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = '="  +1.66%  "'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)
$ws.Range("E4").Formula = '="  +0.35%  "'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("D5").Formula = '="237.37"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = '="  +3.03%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="0.620"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = '="  +0.46%  "'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)
$ws.Range("E7").Formula = '="  +0.28%  "'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("D8").Formula = '="42.20"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Formula = '="  +5.02%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("D9").Formula = '="0.326"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="  +1.01%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("D10").Formula = '="0.0690"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Formula = '="  +1.12%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="0.0987"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = '="  +0.03%  "'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("D12").Formula = '="2.117.27"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="  +1.78%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("B13").Formula = '="WrappedEther"'
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C13").Formula = '="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"'
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$ws.Range("D13").Formula = '="1.854.24"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="  +1.95%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("B14").Formula = '="Chainlink"'
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("C14").Formula = '="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"'
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial(-4163)
$ws.Range("D14").Formula = '="11.36"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Formula = '="  +0.58%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)
$ws.Range("E15").Formula = '="  +0.44%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("D16").Formula = '="4.73"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = '="  +2.62%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("D17").Formula = '="34.957.82"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Formula = '="  +0.36%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)
$ws.Range("D18").Formula = '="69.90"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Formula = '="  +0.33%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("D19").Formula = '="0.0₃0790"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = '="  +0.77%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="240.18"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = '="  -0.20%  "'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="12.13"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Formula = '="  +0.83%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("E22").Formula = '="  +1.06%  "'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("E23").Formula = '="  +0.04%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("E24").Formula = '="  +0.01%  "'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)
$ws.Range("D25").Formula = '="169.51"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = '="  -2.17%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)
$ws.Range("B26").Formula = '="Cosmos"'
$ws.Range("B26").Copy()
$ws.Range("B26").PasteSpecial(-4163)
$ws.Range("C26").Formula = '="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"'
$ws.Range("C26").Copy()
$ws.Range("C26").PasteSpecial(-4163)
$ws.Range("D26").Formula = '="7.96"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Formula = '="  +2.47%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("B27").Formula = '="PancakeSwap"'
$ws.Range("B27").Copy()
$ws.Range("B27").PasteSpecial(-4163)
$ws.Range("C27").Formula = '="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"'
$ws.Range("C27").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("D27").Formula = '="1.83"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = '="  +20.76%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("D28").Formula = '="17.56"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = '="  +1.38%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("D29").Formula = '="0.123"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = '="  -0.14%  "'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("E30").Formula = '="  +0.67%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("D31").Formula = '="0.0552"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Formula = '="  +0.75%  "'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("D32").Formula = '="3.98"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = '="  -0.70%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("E33").Formula = '="  +1.10%  "'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)
$ws.Range("D34").Formula = '="1.69"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Formula = '="  +24.38%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)
$ws.Range("E35").Formula = '="  +8.97%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)
$ws.Range("E36").Formula = '="  +3.70%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("D37").Formula = '="0.777"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = '="  +12.32%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("B38").Formula = '="ARBITRUM"'
$ws.Range("B38").Copy()
$ws.Range("B38").PasteSpecial(-4163)
$ws.Range("C38").Formula = '="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"'
$ws.Range("C38").Copy()
$ws.Range("C38").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="1.07"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Formula = '="  +9.46%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)
$ws.Range("B39").Formula = '="VeChain"'
$ws.Range("B39").Copy()
$ws.Range("B39").PasteSpecial(-4163)
$ws.Range("C39").Formula = '="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"'
$ws.Range("C39").Copy()
$ws.Range("C39").PasteSpecial(-4163)
$ws.Range("D39").Formula = '="0.0202"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Formula = '="  +4.75%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)
$ws.Range("B40").Formula = '="Aave"'
$ws.Range("B40").Copy()
$ws.Range("B40").PasteSpecial(-4163)
$ws.Range("C40").Formula = '="https://coinranking.com/coin/ixgUfzmLR+aave-aave"'
$ws.Range("C40").Copy()
$ws.Range("C40").PasteSpecial(-4163)
$ws.Range("D40").Formula = '="90.04"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = '="  -2.13%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)
$ws.Range("B41").Formula = '="Maker"'
$ws.Range("B41").Copy()
$ws.Range("B41").PasteSpecial(-4163)
$ws.Range("C41").Formula = '="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"'
$ws.Range("C41").Copy()
$ws.Range("C41").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="1.341.00"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = '="  +0.08%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)
$ws.Range("B42").Formula = '="Gas"'
$ws.Range("B42").Copy()
$ws.Range("B42").PasteSpecial(-4163)
$ws.Range("C42").Formula = '="https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"'
$ws.Range("C42").Copy()
$ws.Range("C42").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="13.43"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = '="  +58.44%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)
$ws.Range("D43").Formula = '="14.91"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = '="  +4.33%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)
$ws.Range("D44").Formula = '="2.31"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = '="  +0.38%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)
$ws.Range("D45").Formula = '="2.42"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = '="  +0.08%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)
$ws.Range("B46").Formula = '="MXToken"'
$ws.Range("B46").Copy()
$ws.Range("B46").PasteSpecial(-4163)
$ws.Range("C46").Formula = '="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"'
$ws.Range("C46").Copy()
$ws.Range("C46").PasteSpecial(-4163)
$ws.Range("D46").Formula = '="2.73"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Formula = '="  -1.07%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)
$ws.Range("B47").Formula = '="Kaspa"'
$ws.Range("B47").Copy()
$ws.Range("B47").PasteSpecial(-4163)
$ws.Range("C47").Formula = '="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"'
$ws.Range("C47").Copy()
$ws.Range("C47").PasteSpecial(-4163)
$ws.Range("D47").Formula = '="0.0550"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = '="  +5.62%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)
$ws.Range("D48").Formula = '="6.47"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = '="  +4.40%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)
$ws.Range("D49").Formula = '="2.032.13"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Formula = '="  +1.84%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)
$ws.Range("E50").Formula = '="  +1.77%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)
$ws.Range("E51").Formula = '="  +0.27%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)
